# MorganPatrone2006a C_Stationary generator (alpha non zero) - update experiment values
# for "nuevos experimentos no convexos".
#
# Many of the target cell values are numeric-looking strings (e.g. "0.3",
# "-1", "1.0") that must remain stored as literal TEXT (shared strings),
# matching how the source generator script wrote this workbook. Plain
# `.Value = "0.3"` assignment lets Excel's COM layer auto-coerce such
# strings into real numbers, so we force literal-text semantics by
# stamping the cell as Text (NumberFormat "@") before the assignment and
# then resetting the cell style back to "Normal" immediately afterward,
# so no stray style index lingers on the cell.
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# Sheet tab order (indices are unambiguous - note "Vector_bf" / "Vector_BF"
# differ only by case, and Worksheets.Item(<name>) lookups are
# case-insensitive, so names alone would collide):
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

# --- Restricciones_del_lider ---------------------------------------------
$ws = $wb.Worksheets.Item(2)
Set-TextValue $ws.Range("A2") "2.8499999999999996 - x"
Set-TextValue $ws.Range("B2") "-3.3499999999999996"
Set-TextValue $ws.Range("D2") "0.3"
Set-TextValue $ws.Range("A3") "-2.8499999999999996 + x"
Set-TextValue $ws.Range("B3") "2.3499999999999996"
Set-TextValue $ws.Range("D3") "0.09"

# --- Restricciones_del_follower ------------------------------------------
$ws = $wb.Worksheets.Item(3)
Set-TextValue $ws.Range("A2") "0"
Set-TextValue $ws.Range("B2") "-1"
Set-TextValue $ws.Range("D2") "0.19"
Set-TextValue $ws.Range("E2") "0"
Set-TextValue $ws.Range("F2") "0"
Set-TextValue $ws.Range("A3") "0"
Set-TextValue $ws.Range("B3") "-1"
Set-TextValue $ws.Range("C3") "J_0_LP_v"
Set-TextValue $ws.Range("D3") "0.79"
Set-TextValue $ws.Range("E3") "-8.299999999999999"
Set-TextValue $ws.Range("F3") "-8.5"

# --- Punto_modificado -------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws.Range("A2") "2.8499999999999996"
Set-TextValue $ws.Range("B2") "4.449999999999999"

# --- Vector_bf ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws.Range("A2") "-2.8499999999999996"

# --- Vector_BF ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Range("A2") "0.3699999999999999"
Set-TextValue $ws.Range("A3") "1.0"

# --- Vector_Alpha (numeric cell, not text) -----------------------------------
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 0.8400000000000001
